$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.195.21"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.904.24"
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'306.56"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.5253"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "'0.07251"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'0.8998"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "'0.08382"
$ws.Range("E12").Value = "  +9.98%  "
$ws.Range("D13").Value = "1.903.25"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'94.91"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'5.278"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'0.000008607"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "27.230.60"
$ws.Range("D22").Value = "2.146.84"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'6.434"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'2.285"
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("D27").Value = "'1.754"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("D28").Value = "'18.17"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "'114.84"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "'4.819"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'0.09282"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'0.8097"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").Value = "'0.05065"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "'2.958"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").Value = "'3.387"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").Value = "'2.617"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "'0.5728"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'0.01988"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'6.638"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'8.975"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'117.47"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'0.1512"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.19"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "'37.47"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "'63.80"
$ws.Range("E51").Value = "  +0.32%  "
